$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the calendar dates (every week was off by one day) ---

# Week row 6 (days 1-5): shift everything one column earlier, column J (old day 1) is dropped
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 1
$ws.Range("J6").Clear()

# Week row 8 (days 6-12)
$ws.Range("E8").Value = 12
$ws.Range("F8").Value = 11
$ws.Range("G8").Value = 10
$ws.Range("H8").Value = 9
$ws.Range("I8").Value = 8
$ws.Range("J8").Value = 7
$ws.Range("K8").Value = 6

# Week row 10 (days 13-19)
$ws.Range("E10").Value = 19
$ws.Range("F10").Value = 18
$ws.Range("G10").Value = 17
$ws.Range("H10").Value = 16
$ws.Range("I10").Value = 15
$ws.Range("J10").Value = 14
$ws.Range("K10").Value = 13

# Week row 12 (days 20-26)
$ws.Range("E12").Value = 26
$ws.Range("F12").Value = 25
$ws.Range("G12").Value = 24
$ws.Range("H12").Value = 23
$ws.Range("I12").Value = 22
$ws.Range("J12").Value = 21
$ws.Range("K12").Value = 20

# Week row 14 (days 27-31): gains a new cell G14 for day 31
$ws.Range("G14").Value = 31
$ws.Range("H14").Value = 30
$ws.Range("I14").Value = 20
$ws.Range("J14").Value = 28
$ws.Range("K14").Value = 27

# --- Swap Mirsalari with Bazgir in the first week's recorder assignment row ---
$ws.Range("F7").Value2 = "بازگیر"
$ws.Range("H7").Value2 = "میرسالاری"

# --- Update the selected cell shown when the file is opened ---
$ws.Range("S10").Select()
